# Add a new closing "Authors / Date / paper" slide at the end of the deck,
# and bump the slide id counter the way the authoring session apparently did.

$p = $ppt.ActivePresentation

# --- Insert new slide 19 (blank layout => ppLayoutBlank = 12) ------------
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 12)

# --- Recreate the shape-id numbering (final textbox ends up id=11,
#     "TextBox 10") by creating + discarding 9 throw-away shapes first. ---
for ($i = 1; $i -le 9; $i++) {
    $junk = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $junk.Delete()
}

# --- Add the real textbox with the exact EMU-derived position/size. -----
# (AddTextbox takes points; EMU / 12700 = points)
$left   = 326735   / 12700.0
$top    = 2305615  / 12700.0
$width  = 11538529 / 12700.0
$height = 2246769  / 12700.0

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tf = $tb.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "Authors: Sunoo Park, Albert Kwon, Georg Fuchsbauer, Peter Gaži, Joël Alwen, Krzysztof Pietrzak`rDate: July 2017`rSpaceMint: A Cryptocurrency Based on Proofs of Space.`rhttps://fc18.ifca.ai/preproceedings/78.pdf"
$tr.Font.Size = 28

# --- Paragraph 1: authors list, with surname runs marked bold/plain. -----
$para1 = $tr.Paragraphs(1, 1)
$para1.Characters(1, 7).Font.Bold = $true                 # "Authors"
# ": Sunoo Park, Albert Kwon, Georg " stays regular (already sz=28 from above)

# --- Paragraph 2: "Date: July 2017" -------------------------------------
$para2 = $tr.Paragraphs(2, 1)
$para2.Characters(1, 4).Font.Bold = $true                 # "Date"

# --- Paragraph 3: "SpaceMint: A Cryptocurrency ..." ---------------------
$para3 = $tr.Paragraphs(3, 1)
$para3.Characters(1, 9).Font.Bold = $true                 # "SpaceMint"

# --- Paragraph 4: link, no special formatting ---------------------------

# Shape has no fill.
$tb.Fill.Visible = $false

# Re-affirm the exact target height (AutoSize recomputes it from the
# (simplified) text layout -- Left/Top/Width are untouched by AutoSize
# and already exact from AddTextbox, so leave them alone).
$tb.Height = $height
